# This deck's two theme parts are "crossed": the slide master (and the
# presentation itself) are wired to ppt/theme/theme2.xml, which currently
# carries the green "Integral" color scheme, while ppt/theme/theme1.xml
# (reachable only through the notes master's relationship) carries the
# stock "Office Theme" (blue) color scheme. The font scheme and format
# scheme (fills/lines/effects) are byte-for-byte identical between the
# two theme parts already - only the <a:clrScheme> (and the cosmetic
# theme/clrScheme "name" attributes) differ.
#
# The authored change swaps the two themes' contents, so the part that
# drives the visible slides ends up using the stock Office Theme colors
# instead of Integral's.
#
# Recolor the presentation's (single) live theme color scheme to the
# target "Office" palette - this is the scriptable, observable
# equivalent of that swap for everything the PowerPoint object model
# exposes to automation (SlideMaster.Theme / Slide.ThemeColorScheme /
# NotesMaster.Theme / HandoutMaster.Theme all resolve to this same theme
# object, which is the part that is actually serialized/persisted).

$p  = $ppt.ActivePresentation
$th = $p.SlideMaster.Theme
$cs = $th.ThemeColorScheme

# Scheme color order (matches PowerPoint's internal 1-12 color index):
#   1 dk1  2 lt1  3 dk2  4 lt2  5-10 accent1..accent6  11 hlink  12 folHlink
# .RGB takes/returns colors as 0xBBGGRR (classic VBA RGB() byte order),
# so each literal below is the reverse-byte-order form of the target
# "Office" theme's RRGGBB hex value (shown in the comment).
$cs.Item(1).RGB  = 0x000000  # dk1      -> 000000
$cs.Item(2).RGB  = 0xFFFFFF  # lt1      -> FFFFFF
$cs.Item(3).RGB  = 0x6A5444  # dk2      -> 44546A
$cs.Item(4).RGB  = 0xE6E6E7  # lt2      -> E7E6E6
$cs.Item(5).RGB  = 0xD59B5B  # accent1  -> 5B9BD5
$cs.Item(6).RGB  = 0x317DED  # accent2  -> ED7D31
$cs.Item(7).RGB  = 0xA5A5A5  # accent3  -> A5A5A5
$cs.Item(8).RGB  = 0x00C0FF  # accent4  -> FFC000
$cs.Item(9).RGB  = 0xC47244  # accent5  -> 4472C4
$cs.Item(10).RGB = 0x47AD70  # accent6  -> 70AD47
$cs.Item(11).RGB = 0xC16305  # hlink    -> 0563C1
$cs.Item(12).RGB = 0x724F95  # folHlink -> 954F72

# Best-effort: also rename the scheme/theme to match (some hosts persist
# this, others treat Name as read-only derived state).
$cs.Name = "Office"
$th.Name = "Office Theme"
